$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 810
$ws1.Range("F3").Value = 553
$ws1.Range("F4").Value = 302
$ws1.Range("F6").Value = 1146
$ws1.Range("F7").Value = 327
$ws1.Range("F9").Value = 123
$ws1.Range("F10").Value = 124
$ws1.Range("F11").Value = 1187
$ws1.Range("F14").Value = 859
$ws1.Range("F15").Value = 861
$ws1.Range("F17").Value = 61
$ws1.Range("F18").Value = 73
$ws1.Range("F20").Value = 743
$ws1.Range("F21").Value = 1736
$ws1.Range("F22").Value = 2717
$ws1.Range("F23").Value = 775
$ws1.Range("F25").Value = 2068
$ws1.Range("F27").Value = 2948
$ws1.Range("F28").Value = 561
$ws1.Range("F30").Value = 2
$ws1.Range("F32").Value = 718
$ws1.Range("F33").Value = 142
$ws1.Range("F35").Value = 101
$ws1.Range("F36").Value = 1037
$ws1.Range("F37").Value = 1748
$ws1.Range("F38").Value = 370
$ws1.Range("F40").Value = 547
$ws1.Range("F41").Value = 177
$ws1.Range("F43").Value = 165
$ws1.Range("F44").Value = 36

# Row 26 also flips G26 from a numeric min-price to "Sold out" text
$ws1.Range("F26").Value = 659
$ws1.Range("G26").Value = "已售罄"

# --- Sheet "演出" (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F13").Value = 1

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 810
$ws4.Range("F4").Value = 553
$ws4.Range("F5").Value = 302
$ws4.Range("F7").Value = 1146
$ws4.Range("F8").Value = 327
$ws4.Range("F10").Value = 123
$ws4.Range("F11").Value = 124
$ws4.Range("F12").Value = 1187
$ws4.Range("F14").Value = 859
$ws4.Range("F15").Value = 861
$ws4.Range("F18").Value = 61
$ws4.Range("F20").Value = 73
$ws4.Range("F21").Value = 743
$ws4.Range("F22").Value = 1736
$ws4.Range("F23").Value = 2717
$ws4.Range("F24").Value = 775
$ws4.Range("F28").Value = 2948
$ws4.Range("F29").Value = 561
$ws4.Range("F36").Value = 718
$ws4.Range("F37").Value = 142
$ws4.Range("F39").Value = 101
$ws4.Range("F40").Value = 1037
$ws4.Range("F41").Value = 1748
$ws4.Range("F43").Value = 370
$ws4.Range("F44").Value = 547
$ws4.Range("F45").Value = 177
$ws4.Range("F47").Value = 165
$ws4.Range("F48").Value = 36
